# repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to reflect the repulled data
$ws.Range("F2").Value = -6
$ws.Range("F3").Value = -6
$ws.Range("F4").Value = -7
$ws.Range("F8").Value = -3
$ws.Range("F10").Value = -1
$ws.Range("F11").Value = 4
$ws.Range("F12").Value = 1
$ws.Range("F14").Value = -7
